$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("codes")

$ws.Range("A44").Value = "Társadalmi sokszínűség (etnikai / osztálybeli / más)"
$ws.Range("A45").Value = "Környezeti benyomások (zaj / szmog / más)"
$ws.Range("A48").Value = "Lakossági összetétel / közösségek"
